$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(16, 3).Value = '1002201543'
$ws.Cells.Item(16, 4).Value = 'FRANCISCO JAVIER BENAVIDES SANTOYA'
$ws.Cells.Item(16, 5).Value = '2206'
$ws.Cells.Item(16, 6).Value = 40000
$ws.Cells.Item(16, 7).Value = 1000000

$ws.Cells.Item(17, 3).Value = '1143342434'
$ws.Cells.Item(17, 4).Value = 'JEFFRY ERNEY MORENO CORPUS'
$ws.Cells.Item(17, 5).Value = '2206'
$ws.Cells.Item(17, 6).Value = 40000
$ws.Cells.Item(17, 7).Value = 1000000

$ws.Cells.Item(18, 3).Value = '1051884273'
$ws.Cells.Item(18, 4).Value = 'JOSE FELIX SANDON MIRANDA'
$ws.Cells.Item(18, 5).Value = '2206'
$ws.Cells.Item(18, 6).Value = 40000
$ws.Cells.Item(18, 7).Value = 1000000

$ws.Cells.Item(19, 3).Value = '1001804427'
$ws.Cells.Item(19, 4).Value = 'JORGE LUIS GARCIA DE AVILA'
$ws.Cells.Item(19, 5).Value = '2206'
$ws.Cells.Item(19, 6).Value = 40000
$ws.Cells.Item(19, 7).Value = 1000000

$ws.Cells.Item(20, 3).Value = '1143407879'
$ws.Cells.Item(20, 4).Value = 'MANUEL ANTONIO CANOLES PEREZ'
$ws.Cells.Item(20, 5).Value = '2206'
$ws.Cells.Item(20, 6).Value = 40000
$ws.Cells.Item(20, 7).Value = 1000000

$ws.Cells.Item(21, 3).Value = '1050948392'
$ws.Cells.Item(21, 4).Value = 'PATRICIA DEL ROSARIO CASTELLON PENARANDA'
$ws.Cells.Item(21, 5).Value = '2206'
$ws.Cells.Item(21, 6).Value = 40000
$ws.Cells.Item(21, 7).Value = 1000000

$ws.Cells.Item(22, 3).Value = '1047491722'
$ws.Cells.Item(22, 4).Value = 'ANA LUZ CORONADO ESQUIVIA'
$ws.Cells.Item(22, 5).Value = '2206'
$ws.Cells.Item(22, 6).Value = 48000
$ws.Cells.Item(22, 7).Value = 1200000

$ws.Cells.Item(23, 3).Value = '1002204523'
$ws.Cells.Item(23, 4).Value = 'ALBERTO ELIAS AURAAD ORTEGA'
$ws.Cells.Item(23, 5).Value = '2206'
$ws.Cells.Item(23, 6).Value = 40000
$ws.Cells.Item(23, 7).Value = 1000000

$ws.Cells.Item(24, 3).Value = '91106360'
$ws.Cells.Item(24, 4).Value = 'FRANKLIN EDUARDO VIVAS JEREZ'
$ws.Cells.Item(24, 5).Value = '2206'
$ws.Cells.Item(24, 6).Value = 40000
$ws.Cells.Item(24, 7).Value = 1000000

$ws.Cells.Item(25, 3).Value = '1001944842'
$ws.Cells.Item(25, 4).Value = 'JOHINNER DANIEL MONDUL FONTALVO'
$ws.Cells.Item(25, 5).Value = '2206'
$ws.Cells.Item(25, 6).Value = 40000
$ws.Cells.Item(25, 7).Value = 1000000

$ws.Cells.Item(26, 3).Value = '1007588496'
$ws.Cells.Item(26, 4).Value = 'LUIS FERNANDO YEPES MARTINEZ'
$ws.Cells.Item(26, 5).Value = '2206'
$ws.Cells.Item(26, 6).Value = 40000
$ws.Cells.Item(26, 7).Value = 1000000

$ws.Cells.Item(27, 3).Value = '1002201543'
$ws.Cells.Item(27, 4).Value = 'FRANCISCO JAVIER BENAVIDES SANTOYA'
$ws.Cells.Item(27, 5).Value = '2207'
$ws.Cells.Item(27, 6).Value = 40000
$ws.Cells.Item(27, 7).Value = 1000000

$ws.Cells.Item(28, 3).Value = '1143342434'
$ws.Cells.Item(28, 4).Value = 'JEFFRY ERNEY MORENO CORPUS'
$ws.Cells.Item(28, 5).Value = '2207'
$ws.Cells.Item(28, 6).Value = 40000
$ws.Cells.Item(28, 7).Value = 1000000

$ws.Cells.Item(29, 3).Value = '1051884273'
$ws.Cells.Item(29, 4).Value = 'JOSE FELIX SANDON MIRANDA'
$ws.Cells.Item(29, 5).Value = '2207'
$ws.Cells.Item(29, 6).Value = 40000
$ws.Cells.Item(29, 7).Value = 1000000

$ws.Cells.Item(30, 3).Value = '1001804427'
$ws.Cells.Item(30, 4).Value = 'JORGE LUIS GARCIA DE AVILA'
$ws.Cells.Item(30, 5).Value = '2207'
$ws.Cells.Item(30, 6).Value = 40000
$ws.Cells.Item(30, 7).Value = 1000000

$ws.Cells.Item(31, 3).Value = '1143407879'
$ws.Cells.Item(31, 4).Value = 'MANUEL ANTONIO CANOLES PEREZ'
$ws.Cells.Item(31, 5).Value = '2207'
$ws.Cells.Item(31, 6).Value = 40000
$ws.Cells.Item(31, 7).Value = 1000000

$ws.Cells.Item(32, 3).Value = '1050948392'
$ws.Cells.Item(32, 4).Value = 'PATRICIA DEL ROSARIO CASTELLON PENARANDA'
$ws.Cells.Item(32, 5).Value = '2207'
$ws.Cells.Item(32, 6).Value = 40000
$ws.Cells.Item(32, 7).Value = 1000000

$ws.Cells.Item(33, 3).Value = '1002204523'
$ws.Cells.Item(33, 4).Value = 'ALBERTO ELIAS AURAAD ORTEGA'
$ws.Cells.Item(33, 5).Value = '2207'
$ws.Cells.Item(33, 6).Value = 40000
$ws.Cells.Item(33, 7).Value = 1000000

$ws.Cells.Item(34, 3).Value = '91106360'
$ws.Cells.Item(34, 4).Value = 'FRANKLIN EDUARDO VIVAS JEREZ'
$ws.Cells.Item(34, 5).Value = '2207'
$ws.Cells.Item(34, 6).Value = 40000
$ws.Cells.Item(34, 7).Value = 1000000

$ws.Cells.Item(35, 3).Value = '1001944842'
$ws.Cells.Item(35, 4).Value = 'JOHINNER DANIEL MONDUL FONTALVO'
$ws.Cells.Item(35, 5).Value = '2207'
$ws.Cells.Item(35, 6).Value = 40000
$ws.Cells.Item(35, 7).Value = 1000000

$ws.Cells.Item(36, 3).Value = '1002201543'
$ws.Cells.Item(36, 4).Value = 'FRANCISCO JAVIER BENAVIDES SANTOYA'
$ws.Cells.Item(36, 5).Value = '2208'
$ws.Cells.Item(36, 6).Value = 40000
$ws.Cells.Item(36, 7).Value = 1000000

$ws.Cells.Item(37, 3).Value = '1143342434'
$ws.Cells.Item(37, 4).Value = 'JEFFRY ERNEY MORENO CORPUS'
$ws.Cells.Item(37, 5).Value = '2208'
$ws.Cells.Item(37, 6).Value = 40000
$ws.Cells.Item(37, 7).Value = 1000000

$ws.Cells.Item(38, 3).Value = '1051884273'
$ws.Cells.Item(38, 4).Value = 'JOSE FELIX SANDON MIRANDA'
$ws.Cells.Item(38, 5).Value = '2208'
$ws.Cells.Item(38, 6).Value = 40000
$ws.Cells.Item(38, 7).Value = 1000000

$ws.Cells.Item(39, 3).Value = '1001804427'
$ws.Cells.Item(39, 4).Value = 'JORGE LUIS GARCIA DE AVILA'
$ws.Cells.Item(39, 5).Value = '2208'
$ws.Cells.Item(39, 6).Value = 40000
$ws.Cells.Item(39, 7).Value = 1000000

$ws.Cells.Item(40, 3).Value = '1143407879'
$ws.Cells.Item(40, 4).Value = 'MANUEL ANTONIO CANOLES PEREZ'
$ws.Cells.Item(40, 5).Value = '2208'
$ws.Cells.Item(40, 6).Value = 40000
$ws.Cells.Item(40, 7).Value = 1000000

$ws.Cells.Item(41, 3).Value = '1050948392'
$ws.Cells.Item(41, 4).Value = 'PATRICIA DEL ROSARIO CASTELLON PENARANDA'
$ws.Cells.Item(41, 5).Value = '2208'
$ws.Cells.Item(41, 6).Value = 40000
$ws.Cells.Item(41, 7).Value = 1000000

$ws.Cells.Item(42, 3).Value = '1002204523'
$ws.Cells.Item(42, 4).Value = 'ALBERTO ELIAS AURAAD ORTEGA'
$ws.Cells.Item(42, 5).Value = '2208'
$ws.Cells.Item(42, 6).Value = 40000
$ws.Cells.Item(42, 7).Value = 1000000

$ws.Cells.Item(43, 3).Value = '91106360'
$ws.Cells.Item(43, 4).Value = 'FRANKLIN EDUARDO VIVAS JEREZ'
$ws.Cells.Item(43, 5).Value = '2208'
$ws.Cells.Item(43, 6).Value = 40000
$ws.Cells.Item(43, 7).Value = 1000000

$ws.Cells.Item(44, 3).Value = '1001944842'
$ws.Cells.Item(44, 4).Value = 'JOHINNER DANIEL MONDUL FONTALVO'
$ws.Cells.Item(44, 5).Value = '2208'
$ws.Cells.Item(44, 6).Value = 40000
$ws.Cells.Item(44, 7).Value = 1000000

$ws.Cells.Item(45, 3).Value = '1002201543'
$ws.Cells.Item(45, 4).Value = 'FRANCISCO JAVIER BENAVIDES SANTOYA'
$ws.Cells.Item(45, 5).Value = '2209'
$ws.Cells.Item(45, 6).Value = 40000
$ws.Cells.Item(45, 7).Value = 1000000

$ws.Cells.Item(46, 3).Value = '1143342434'
$ws.Cells.Item(46, 4).Value = 'JEFFRY ERNEY MORENO CORPUS'
$ws.Cells.Item(46, 5).Value = '2209'
$ws.Cells.Item(46, 6).Value = 40000
$ws.Cells.Item(46, 7).Value = 1000000

$ws.Cells.Item(47, 3).Value = '1051884273'
$ws.Cells.Item(47, 4).Value = 'JOSE FELIX SANDON MIRANDA'
$ws.Cells.Item(47, 5).Value = '2209'
$ws.Cells.Item(47, 6).Value = 40000
$ws.Cells.Item(47, 7).Value = 1000000

$ws.Cells.Item(48, 3).Value = '1001804427'
$ws.Cells.Item(48, 4).Value = 'JORGE LUIS GARCIA DE AVILA'
$ws.Cells.Item(48, 5).Value = '2209'
$ws.Cells.Item(48, 6).Value = 40000
$ws.Cells.Item(48, 7).Value = 1000000

$ws.Cells.Item(49, 3).Value = '1143407879'
$ws.Cells.Item(49, 4).Value = 'MANUEL ANTONIO CANOLES PEREZ'
$ws.Cells.Item(49, 5).Value = '2209'
$ws.Cells.Item(49, 6).Value = 40000
$ws.Cells.Item(49, 7).Value = 1000000

$ws.Cells.Item(50, 3).Value = '1050948392'
$ws.Cells.Item(50, 4).Value = 'PATRICIA DEL ROSARIO CASTELLON PENARANDA'
$ws.Cells.Item(50, 5).Value = '2209'
$ws.Cells.Item(50, 6).Value = 40000
$ws.Cells.Item(50, 7).Value = 1000000

$ws.Cells.Item(51, 3).Value = '1002204523'
$ws.Cells.Item(51, 4).Value = 'ALBERTO ELIAS AURAAD ORTEGA'
$ws.Cells.Item(51, 5).Value = '2209'
$ws.Cells.Item(51, 6).Value = 40000
$ws.Cells.Item(51, 7).Value = 1000000

$ws.Cells.Item(52, 3).Value = '91106360'
$ws.Cells.Item(52, 4).Value = 'FRANKLIN EDUARDO VIVAS JEREZ'
$ws.Cells.Item(52, 5).Value = '2209'
$ws.Cells.Item(52, 6).Value = 40000
$ws.Cells.Item(52, 7).Value = 1000000

$ws.Cells.Item(53, 3).Value = '1001944842'
$ws.Cells.Item(53, 4).Value = 'JOHINNER DANIEL MONDUL FONTALVO'
$ws.Cells.Item(53, 5).Value = '2209'
$ws.Cells.Item(53, 6).Value = 40000
$ws.Cells.Item(53, 7).Value = 1000000

$ws.Cells.Item(54, 3).Value = '1002201543'
$ws.Cells.Item(54, 4).Value = 'FRANCISCO JAVIER BENAVIDES SANTOYA'
$ws.Cells.Item(54, 5).Value = '2210'
$ws.Cells.Item(54, 6).Value = 40000
$ws.Cells.Item(54, 7).Value = 1000000

$ws.Cells.Item(55, 3).Value = '1143342434'
$ws.Cells.Item(55, 4).Value = 'JEFFRY ERNEY MORENO CORPUS'
$ws.Cells.Item(55, 5).Value = '2210'
$ws.Cells.Item(55, 6).Value = 40000
$ws.Cells.Item(55, 7).Value = 1000000

$ws.Cells.Item(56, 3).Value = '1051884273'
$ws.Cells.Item(56, 4).Value = 'JOSE FELIX SANDON MIRANDA'
$ws.Cells.Item(56, 5).Value = '2210'
$ws.Cells.Item(56, 6).Value = 40000
$ws.Cells.Item(56, 7).Value = 1000000

$ws.Cells.Item(57, 3).Value = '1001804427'
$ws.Cells.Item(57, 4).Value = 'JORGE LUIS GARCIA DE AVILA'
$ws.Cells.Item(57, 5).Value = '2210'
$ws.Cells.Item(57, 6).Value = 40000
$ws.Cells.Item(57, 7).Value = 1000000

$ws.Cells.Item(58, 3).Value = '1143407879'
$ws.Cells.Item(58, 4).Value = 'MANUEL ANTONIO CANOLES PEREZ'
$ws.Cells.Item(58, 5).Value = '2210'
$ws.Cells.Item(58, 6).Value = 40000
$ws.Cells.Item(58, 7).Value = 1000000

$ws.Cells.Item(59, 3).Value = '1050948392'
$ws.Cells.Item(59, 4).Value = 'PATRICIA DEL ROSARIO CASTELLON PENARANDA'
$ws.Cells.Item(59, 5).Value = '2210'
$ws.Cells.Item(59, 6).Value = 40000
$ws.Cells.Item(59, 7).Value = 1000000

$ws.Cells.Item(60, 3).Value = '1002204523'
$ws.Cells.Item(60, 4).Value = 'ALBERTO ELIAS AURAAD ORTEGA'
$ws.Cells.Item(60, 5).Value = '2210'
$ws.Cells.Item(60, 6).Value = 40000
$ws.Cells.Item(60, 7).Value = 1000000

$ws.Cells.Item(61, 3).Value = '91106360'
$ws.Cells.Item(61, 4).Value = 'FRANKLIN EDUARDO VIVAS JEREZ'
$ws.Cells.Item(61, 5).Value = '2210'
$ws.Cells.Item(61, 6).Value = 40000
$ws.Cells.Item(61, 7).Value = 1000000

$ws.Cells.Item(62, 3).Value = '1001944842'
$ws.Cells.Item(62, 4).Value = 'JOHINNER DANIEL MONDUL FONTALVO'
$ws.Cells.Item(62, 5).Value = '2210'
$ws.Cells.Item(62, 6).Value = 40000
$ws.Cells.Item(62, 7).Value = 1000000

$ws.Cells.Item(63, 3).Value = '1002201543'
$ws.Cells.Item(63, 4).Value = 'FRANCISCO JAVIER BENAVIDES SANTOYA'
$ws.Cells.Item(63, 5).Value = '2211'
$ws.Cells.Item(63, 6).Value = 40000
$ws.Cells.Item(63, 7).Value = 1000000

$ws.Cells.Item(64, 3).Value = '1143342434'
$ws.Cells.Item(64, 4).Value = 'JEFFRY ERNEY MORENO CORPUS'
$ws.Cells.Item(64, 5).Value = '2211'
$ws.Cells.Item(64, 6).Value = 40000
$ws.Cells.Item(64, 7).Value = 1000000

$ws.Cells.Item(65, 3).Value = '1051884273'
$ws.Cells.Item(65, 4).Value = 'JOSE FELIX SANDON MIRANDA'
$ws.Cells.Item(65, 5).Value = '2211'
$ws.Cells.Item(65, 6).Value = 40000
$ws.Cells.Item(65, 7).Value = 1000000

$ws.Cells.Item(66, 3).Value = '1001804427'
$ws.Cells.Item(66, 4).Value = 'JORGE LUIS GARCIA DE AVILA'
$ws.Cells.Item(66, 5).Value = '2211'
$ws.Cells.Item(66, 6).Value = 40000
$ws.Cells.Item(66, 7).Value = 1000000

$ws.Cells.Item(67, 3).Value = '1143407879'
$ws.Cells.Item(67, 4).Value = 'MANUEL ANTONIO CANOLES PEREZ'
$ws.Cells.Item(67, 5).Value = '2211'
$ws.Cells.Item(67, 6).Value = 40000
$ws.Cells.Item(67, 7).Value = 1000000

$ws.Cells.Item(68, 3).Value = '1050948392'
$ws.Cells.Item(68, 4).Value = 'PATRICIA DEL ROSARIO CASTELLON PENARANDA'
$ws.Cells.Item(68, 5).Value = '2211'
$ws.Cells.Item(68, 6).Value = 40000
$ws.Cells.Item(68, 7).Value = 1000000

$ws.Cells.Item(69, 3).Value = '1002204523'
$ws.Cells.Item(69, 4).Value = 'ALBERTO ELIAS AURAAD ORTEGA'
$ws.Cells.Item(69, 5).Value = '2211'
$ws.Cells.Item(69, 6).Value = 40000
$ws.Cells.Item(69, 7).Value = 1000000

$ws.Cells.Item(70, 3).Value = '91106360'
$ws.Cells.Item(70, 4).Value = 'FRANKLIN EDUARDO VIVAS JEREZ'
$ws.Cells.Item(70, 5).Value = '2211'
$ws.Cells.Item(70, 6).Value = 40000
$ws.Cells.Item(70, 7).Value = 1000000

$ws.Cells.Item(71, 3).Value = '1001944842'
$ws.Cells.Item(71, 4).Value = 'JOHINNER DANIEL MONDUL FONTALVO'
$ws.Cells.Item(71, 5).Value = '2211'
$ws.Cells.Item(71, 6).Value = 40000
$ws.Cells.Item(71, 7).Value = 1000000

$ws.Cells.Item(72, 3).Value = '1002201543'
$ws.Cells.Item(72, 4).Value = 'FRANCISCO JAVIER BENAVIDES SANTOYA'
$ws.Cells.Item(72, 5).Value = '2212'
$ws.Cells.Item(72, 6).Value = 40000
$ws.Cells.Item(72, 7).Value = 1000000

$ws.Cells.Item(73, 3).Value = '1143342434'
$ws.Cells.Item(73, 4).Value = 'JEFFRY ERNEY MORENO CORPUS'
$ws.Cells.Item(73, 5).Value = '2212'
$ws.Cells.Item(73, 6).Value = 40000
$ws.Cells.Item(73, 7).Value = 1000000

$ws.Cells.Item(74, 3).Value = '1051884273'
$ws.Cells.Item(74, 4).Value = 'JOSE FELIX SANDON MIRANDA'
$ws.Cells.Item(74, 5).Value = '2212'
$ws.Cells.Item(74, 6).Value = 40000
$ws.Cells.Item(74, 7).Value = 1000000

$ws.Cells.Item(75, 3).Value = '1001804427'
$ws.Cells.Item(75, 4).Value = 'JORGE LUIS GARCIA DE AVILA'
$ws.Cells.Item(75, 5).Value = '2212'
$ws.Cells.Item(75, 6).Value = 40000
$ws.Cells.Item(75, 7).Value = 1000000

$ws.Cells.Item(76, 3).Value = '1143407879'
$ws.Cells.Item(76, 4).Value = 'MANUEL ANTONIO CANOLES PEREZ'
$ws.Cells.Item(76, 5).Value = '2212'
$ws.Cells.Item(76, 6).Value = 40000
$ws.Cells.Item(76, 7).Value = 1000000

$ws.Cells.Item(77, 3).Value = '1050948392'
$ws.Cells.Item(77, 4).Value = 'PATRICIA DEL ROSARIO CASTELLON PENARANDA'
$ws.Cells.Item(77, 5).Value = '2212'
$ws.Cells.Item(77, 6).Value = 40000
$ws.Cells.Item(77, 7).Value = 1000000

$ws.Cells.Item(78, 3).Value = '1002204523'
$ws.Cells.Item(78, 4).Value = 'ALBERTO ELIAS AURAAD ORTEGA'
$ws.Cells.Item(78, 5).Value = '2212'
$ws.Cells.Item(78, 6).Value = 40000
$ws.Cells.Item(78, 7).Value = 1000000

$ws.Cells.Item(79, 3).Value = '91106360'
$ws.Cells.Item(79, 4).Value = 'FRANKLIN EDUARDO VIVAS JEREZ'
$ws.Cells.Item(79, 5).Value = '2212'
$ws.Cells.Item(79, 6).Value = 40000
$ws.Cells.Item(79, 7).Value = 1000000

$ws.Cells.Item(80, 3).Value = '1001944842'
$ws.Cells.Item(80, 4).Value = 'JOHINNER DANIEL MONDUL FONTALVO'
$ws.Cells.Item(80, 5).Value = '2212'
$ws.Cells.Item(80, 6).Value = 40000
$ws.Cells.Item(80, 7).Value = 1000000

$ws.Cells.Item(81, 3).Value = '1002201543'
$ws.Cells.Item(81, 4).Value = 'FRANCISCO JAVIER BENAVIDES SANTOYA'
$ws.Cells.Item(81, 5).Value = '2301'
$ws.Cells.Item(81, 6).Value = 33333
$ws.Cells.Item(81, 7).Value = 1000000

$ws.Cells.Item(82, 3).Value = '1143342434'
$ws.Cells.Item(82, 4).Value = 'JEFFRY ERNEY MORENO CORPUS'
$ws.Cells.Item(82, 5).Value = '2301'
$ws.Cells.Item(82, 6).Value = 33333
$ws.Cells.Item(82, 7).Value = 1000000

$ws.Cells.Item(83, 3).Value = '1051884273'
$ws.Cells.Item(83, 4).Value = 'JOSE FELIX SANDON MIRANDA'
$ws.Cells.Item(83, 5).Value = '2301'
$ws.Cells.Item(83, 6).Value = 33333
$ws.Cells.Item(83, 7).Value = 1000000

$ws.Cells.Item(84, 3).Value = '1001804427'
$ws.Cells.Item(84, 4).Value = 'JORGE LUIS GARCIA DE AVILA'
$ws.Cells.Item(84, 5).Value = '2301'
$ws.Cells.Item(84, 6).Value = 33333
$ws.Cells.Item(84, 7).Value = 1000000

$ws.Cells.Item(85, 3).Value = '1143407879'
$ws.Cells.Item(85, 4).Value = 'MANUEL ANTONIO CANOLES PEREZ'
$ws.Cells.Item(85, 5).Value = '2301'
$ws.Cells.Item(85, 6).Value = 33333
$ws.Cells.Item(85, 7).Value = 1000000

$ws.Cells.Item(86, 3).Value = '1050948392'
$ws.Cells.Item(86, 4).Value = 'PATRICIA DEL ROSARIO CASTELLON PENARANDA'
$ws.Cells.Item(86, 5).Value = '2301'
$ws.Cells.Item(86, 6).Value = 33333
$ws.Cells.Item(86, 7).Value = 1000000

$ws.Cells.Item(87, 3).Value = '1002204523'
$ws.Cells.Item(87, 4).Value = 'ALBERTO ELIAS AURAAD ORTEGA'
$ws.Cells.Item(87, 5).Value = '2301'
$ws.Cells.Item(87, 6).Value = 33333
$ws.Cells.Item(87, 7).Value = 1000000

$ws.Cells.Item(88, 3).Value = '91106360'
$ws.Cells.Item(88, 4).Value = 'FRANKLIN EDUARDO VIVAS JEREZ'
$ws.Cells.Item(88, 5).Value = '2301'
$ws.Cells.Item(88, 6).Value = 33333
$ws.Cells.Item(88, 7).Value = 1000000

$ws.Cells.Item(89, 3).Value = '1001944842'
$ws.Cells.Item(89, 4).Value = 'JOHINNER DANIEL MONDUL FONTALVO'
$ws.Cells.Item(89, 5).Value = '2301'
$ws.Cells.Item(89, 6).Value = 33333
$ws.Cells.Item(89, 7).Value = 1000000
